$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.212.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.843.93'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.67'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6713'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07426'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.35%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.86'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07720'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.831.79'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.007'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6742'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '86.11'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.121'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008318'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.159.12'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.75'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.192'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.94%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.48'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1403'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.01'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.510'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.179'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.068'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.189'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05311'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.89%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.878'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.18%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7574'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.47%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.328.05'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.67%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9205'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.27%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.35'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.07976'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +14.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.973.60'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5162'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.774'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.99%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.86'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05946'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.38%  '
